# Apply "Nieuws items stylen + DB connection" update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the next time-tracking entry (row 10), matching the
# formatting already used by the rows above it
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A10").Value = 42394

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Nieuws items stylen + mysql DB onlin zetten"

# Recalculate so the Totaal formula in D2 picks up the new hours
$excel.Calculate()

# Move the view / selection the way the author left it
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C16").Select()
